$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Bugs")

# Update state values: "New" -> "Review" / "InProgress" for specific rows
$ws.Range("C2").Value = "Review"
$ws.Range("C3").Value = "Review"
$ws.Range("C4").Value = "InProgress"
$ws.Range("C7").Value = "Review"

# Add new comments/descriptions in column F
$ws.Range("F2").Value = "Localization file added in common"
$ws.Range("F7").Value = "Added button at the bottom of entity page"

# Update selection on the Bugs sheet
$ws.Range("C3").Select()
